$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 12 (new local extremum for Altai - Choysky district) by
# duplicating the formatting of the last existing row (11) and then
# overwriting with the new record's values.
$ws.Rows.Item(11).Copy()
$ws.Rows.Item(12).Insert(-4121)  # xlShiftDown
$excel.CutCopyMode = $false

# Column G in row 11 uses a plain-integer format (matches style used by the
# "retailturnover" column elsewhere); the new row instead needs the same
# "0.000" style used by the rest of the numeric columns, so copy that
# format in from a neighbouring cell before writing the value.
$ws.Range("F11").Copy()
$ws.Range("G12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A12").Value = 84645000
$ws.Range("B12").Value = "Чойский МР"
$ws.Range("C12").Value = 2018
$ws.Range("D12").Value = -136
$ws.Range("E12").Value = 8315
$ws.Range("F12").Value = 0.17714972940469031
$ws.Range("G12").Value = 20566.364119999998
$ws.Range("H12").Value = 0.69592303066746852
$ws.Range("I12").Value = [double]"1.7919422730006015E-2"
$ws.Range("J12").Value = 49.556223692122671
$ws.Range("K12").Value = 19.8
$ws.Range("L12").Value = [double]"2.1647624774503911E-3"
$ws.Range("M12").Value = [double]"8.4185207456404084E-4"
$ws.Range("N12").Value = [double]"1.8989777510523151E-2"
$ws.Range("O12").Value = 1.2476247745039086
$ws.Range("P12").Value = 0.56698376428141917
$ws.Range("Q12").Value = 22.949180844257366
$ws.Range("R12").Value = [double]"1.2026458208057728E-3"
$ws.Range("S12").Value = [double]"7.4323511725796751E-2"
$ws.Range("T12").Value = 29.9954811785929

# Column G width adjustment (engine quantises ColumnWidth to 1/6-character
# steps, so this is the closest achievable match to the target 10.7109375
# stored width)
$ws.Columns.Item(7).ColumnWidth = 9.86

# Update selection to match the recorded cursor position after edit
$ws.Range("G15").Select()
